$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 683.4286
$ws.Range("I28").Value = 202.75
$ws.Range("K28").Value = 202.75
$ws.Range("M28").Value = 282.25
$ws.Range("H69").Value = 1498.1034
$ws.Range("I69").Value = 1650
$ws.Range("J69").Value = 1486.8518
$ws.Range("K69").Value = 4950
$ws.Range("L69").Value = 4460.555399999999
$ws.Range("M69").Value = -4076
$ws.Range("N69").Value = -6208.555399999999
$ws.Range("H72").Value = 1498.1034
$ws.Range("I72").Value = 1650
$ws.Range("J72").Value = 1486.8518
$ws.Range("K72").Value = 14850
$ws.Range("L72").Value = 13381.6662
$ws.Range("M72").Value = -10482
$ws.Range("N72").Value = -22117.6662
$ws.Range("H86").Value = 10636.0625
$ws.Range("I86").Value = 5338.6
$ws.Range("K86").Value = 5338.6
$ws.Range("M86").Value = -4215.6
$ws.Range("H89").Value = 10636.0625
$ws.Range("I89").Value = 5338.6
$ws.Range("K89").Value = 26693
$ws.Range("M89").Value = -21077
$ws.Range("H98").Value = 2825
$ws.Range("I98").Value = 2000
$ws.Range("J98").Value = 3100
$ws.Range("K98").Value = 2000
$ws.Range("L98").Value = 3100
$ws.Range("M98").Value = -502
$ws.Range("N98").Value = -6096
$ws.Range("H106").Value = 11496566
$ws.Range("I106").Value = 19609222
$ws.Range("J106").Value = 3636
$ws.Range("K106").Value = 19609222
$ws.Range("L106").Value = 3636
$ws.Range("M106").Value = -19608591
$ws.Range("N106").Value = -4898
$ws.Range("H111").Value = 2570.375
$ws.Range("I111").Value = 2855.923
$ws.Range("J111").Value = 1333
$ws.Range("K111").Value = 8567.769
$ws.Range("L111").Value = 3999
$ws.Range("M111").Value = -5500.769
$ws.Range("N111").Value = -10133
$ws.Range("H122").Value = 2825
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -14200
$ws.Range("H129").Value = 501343.6
$ws.Range("J129").Value = 626605.1
$ws.Range("L129").Value = 1879815.3
$ws.Range("N129").Value = -1889815.3
$ws.Range("H132").Value = 2649.4688
$ws.Range("I132").Value = 2649.4688
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7948.4064
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -5418.4064
$ws.Range("H133").Value = 47736.332
$ws.Range("J133").Value = 47736.332
$ws.Range("L133").Value = 47736.332
$ws.Range("N133").Value = -57856.332
$ws.Range("H138").Value = 2124.0488
$ws.Range("I138").Value = 1766.4445
$ws.Range("J138").Value = 2224.625
$ws.Range("K138").Value = 5299.333500000001
$ws.Range("L138").Value = 6673.875
$ws.Range("M138").Value = -159.3335000000006
$ws.Range("N138").Value = -16953.875
$ws.Range("H139").Value = 49933.332
$ws.Range("J139").Value = 49933.332
$ws.Range("L139").Value = 49933.332
$ws.Range("N139").Value = -60213.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5101.559
$ws.Range("I32").Value = 5000.9434
$ws.Range("J32").Value = 5990.3335
$ws.Range("K32").Value = 5000.9434
$ws.Range("L32").Value = 5990.3335
$ws.Range("M32").Value = -4713.9434
$ws.Range("N32").Value = -6564.3335
$ws.Range("H61").Value = 1419.5312
$ws.Range("I61").Value = 1257.5
$ws.Range("J61").Value = 3850
$ws.Range("K61").Value = 1257.5
$ws.Range("L61").Value = 3850
$ws.Range("M61").Value = -1045.5
$ws.Range("N61").Value = -4274
$ws.Range("H122").Value = 1458.2444
$ws.Range("I122").Value = 1155.9117
$ws.Range("J122").Value = 2392.7273
$ws.Range("K122").Value = 3467.7351
$ws.Range("L122").Value = 7178.1819
$ws.Range("M122").Value = -1017.7351
$ws.Range("N122").Value = -12078.1819
$ws.Range("H132").Value = 26769.883
$ws.Range("I132").Value = 1516.25
$ws.Range("J132").Value = 127784.414
$ws.Range("K132").Value = 4548.75
$ws.Range("L132").Value = 383353.242
$ws.Range("M132").Value = -2018.75
$ws.Range("N132").Value = -388413.242
$ws.Range("H136").Value = 1419.5312
$ws.Range("I136").Value = 1257.5
$ws.Range("J136").Value = 3850
$ws.Range("K136").Value = 3772.5
$ws.Range("L136").Value = 11550
$ws.Range("M136").Value = -1222.5
$ws.Range("N136").Value = -16650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 761.9706
$ws.Range("I107").Value = 582.619
$ws.Range("J107").Value = 1051.6923
$ws.Range("K107").Value = 582.619
$ws.Range("L107").Value = 1051.6923
$ws.Range("M107").Value = 1337.381
$ws.Range("N107").Value = -4891.6923
$ws.Range("H134").Value = 7193.6875
$ws.Range("I134").Value = 8045
$ws.Range("J134").Value = 3504.6667
$ws.Range("K134").Value = 24135
$ws.Range("L134").Value = 10514.0001
$ws.Range("M134").Value = -21600
$ws.Range("N134").Value = -15584.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11242.069
$ws.Range("I31").Value = 13878.839
$ws.Range("J31").Value = 4430.4165
$ws.Range("K31").Value = 13878.839
$ws.Range("L31").Value = 4430.4165
$ws.Range("M31").Value = -13583.839
$ws.Range("N31").Value = -5020.4165
$ws.Range("H34").Value = 11242.069
$ws.Range("I34").Value = 13878.839
$ws.Range("J34").Value = 4430.4165
$ws.Range("K34").Value = 13878.839
$ws.Range("L34").Value = 4430.4165
$ws.Range("M34").Value = -13676.839
$ws.Range("N34").Value = -4834.4165
$ws.Range("H58").Value = 18900.393
$ws.Range("I58").Value = 1233.25
$ws.Range("J58").Value = 63068.25
$ws.Range("K58").Value = 1233.25
$ws.Range("L58").Value = 63068.25
$ws.Range("M58").Value = -1030.25
$ws.Range("N58").Value = -63474.25
$ws.Range("H74").Value = 25573.166
$ws.Range("J74").Value = 25573.166
$ws.Range("L74").Value = 25573.166
$ws.Range("N74").Value = -27321.166
$ws.Range("H77").Value = 25573.166
$ws.Range("J77").Value = 25573.166
$ws.Range("L77").Value = 76719.49800000001
$ws.Range("N77").Value = -85455.49800000001
$ws.Range("H94").Value = 2582.0715
$ws.Range("I94").Value = 733.6667
$ws.Range("K94").Value = 733.6667
$ws.Range("M94").Value = -282.6667
$ws.Range("H122").Value = 1097.5625
$ws.Range("J122").Value = 1285.3572
$ws.Range("L122").Value = 3856.0716
$ws.Range("N122").Value = -8756.071599999999
$ws.Range("H132").Value = 15424.895
$ws.Range("I132").Value = 16640
$ws.Range("J132").Value = 7405.2
$ws.Range("K132").Value = 49920
$ws.Range("L132").Value = 22215.6
$ws.Range("M132").Value = -47390
$ws.Range("N132").Value = -27275.6
$ws.Range("H134").Value = 672.64105
$ws.Range("I134").Value = 566.2258
$ws.Range("J134").Value = 1085
$ws.Range("K134").Value = 1698.6774
$ws.Range("L134").Value = 3255
$ws.Range("M134").Value = 836.3226
$ws.Range("N134").Value = -8325
$ws.Range("H136").Value = 18900.393
$ws.Range("I136").Value = 1233.25
$ws.Range("J136").Value = 63068.25
$ws.Range("K136").Value = 3699.75
$ws.Range("L136").Value = 189204.75
$ws.Range("M136").Value = -1149.75
$ws.Range("N136").Value = -194304.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1666.7778
$ws.Range("J23").Value = 1687.625
$ws.Range("L23").Value = 5062.875
$ws.Range("N23").Value = -5532.875
$ws.Range("H44").Value = 329
$ws.Range("I44").Value = 329
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 987
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -589
$ws.Range("H86").Value = 55555804
$ws.Range("I86").Value = 262.25
$ws.Range("J86").Value = 100000240
$ws.Range("K86").Value = 786.75
$ws.Range("L86").Value = 300000720
$ws.Range("M86").Value = 399.25
$ws.Range("N86").Value = -300003092
$ws.Range("H89").Value = 55555804
$ws.Range("I89").Value = 262.25
$ws.Range("J89").Value = 100000240
$ws.Range("K89").Value = 2360.25
$ws.Range("L89").Value = 900002160
$ws.Range("M89").Value = 3567.75
$ws.Range("N89").Value = -900014016
$ws.Range("H131").Value = 143671.73
$ws.Range("J131").Value = 152349.86
$ws.Range("L131").Value = 457049.58
$ws.Range("N131").Value = -467129.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5128443.5
$ws.Range("I107").Value = 266.54544
$ws.Range("K107").Value = 266.54544
$ws.Range("M107").Value = 1653.45456
$ws.Range("H122").Value = 83334710
$ws.Range("I122").Value = 37037880
$ws.Range("J122").Value = 142859220
$ws.Range("K122").Value = 111113640
$ws.Range("L122").Value = 428577660
$ws.Range("M122").Value = -111111190
$ws.Range("N122").Value = -428582560
$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 389.66666
$ws.Range("I16").Value = 437.16666
$ws.Range("K16").Value = 437.16666
$ws.Range("M16").Value = -267.16666
$ws.Range("H40").Value = 3252.25
$ws.Range("I40").Value = 2603.1428
$ws.Range("J40").Value = 4161
$ws.Range("K40").Value = 2603.1428
$ws.Range("L40").Value = 4161
$ws.Range("M40").Value = -2467.1428
$ws.Range("N40").Value = -4433
$ws.Range("H132").Value = 1888.8462
$ws.Range("I132").Value = 1596
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 4788
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -2258
$ws.Range("N132").Value = -15558.5
$ws.Range("H136").Value = 14718.919
$ws.Range("I136").Value = 21816.25
$ws.Range("J136").Value = 1616.1538
$ws.Range("K136").Value = 65448.75
$ws.Range("L136").Value = 4848.4614
$ws.Range("M136").Value = -62898.75
$ws.Range("N136").Value = -9948.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 100000790
$ws.Range("I81").Value = 847.4286
$ws.Range("J81").Value = 333334000
$ws.Range("K81").Value = 1694.8572
$ws.Range("L81").Value = 666668000
$ws.Range("M81").Value = -633.8571999999999
$ws.Range("N81").Value = -666670122
$ws.Range("H84").Value = 100000790
$ws.Range("I84").Value = 847.4286
$ws.Range("J84").Value = 333334000
$ws.Range("K84").Value = 8474.286
$ws.Range("L84").Value = 3333340000
$ws.Range("M84").Value = -3170.286
$ws.Range("N84").Value = -3333350608
$ws.Range("H122").Value = 1317.8422
$ws.Range("J122").Value = 1499.9
$ws.Range("L122").Value = 4499.700000000001
$ws.Range("N122").Value = -9399.700000000001
